$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("T5").ClearContents()

$ws.Range("B2").Value = "'4"
$ws.Range("D2").Value = 0.0788
$ws.Range("E2").Value = 0.116
$ws.Range("F2").Value = -0.10725
$ws.Range("G2").Value = 0.1507861242825056
$ws.Range("H2").Value = 0.1507861242825056
$ws.Range("I2").Value = 0.1407976041926629
$ws.Range("J2").Value = 0.1202715213990252
$ws.Range("K2").Value = 209.728
$ws.Range("L2").Value = 0.1046808085849763
$ws.Range("M2").Value = 56.8
$ws.Range("N2").Value = 0.03263598807177619
$ws.Range("O2").Value = 0.2708269758925846
$ws.Range("P2").Value = 56.8
$ws.Range("Q2").Value = 0.03263598807177619
$ws.Range("R2").Value = 0.2708269758925846
$ws.Range("U2").Value = 365.82
$ws.Range("V2").Value = 0.2101918513453726
$ws.Range("W2").Value = 0.2078169457093224
$ws.Range("X2").Value = 0.04848642154645767
$ws.Range("Y2").Value = 0.1593305241628648
$ws.Range("Z2").Value = 2.256755341137828
$ws.Range("AA2").Value = 0.1437635940206736
$ws.Range("AB2").Value = 0.04823145852529305
$ws.Range("AC2").Value = 0.0955222249095907
$ws.Range("AD2").Value = 12.554
$ws.Range("AE2").Value = 0
$ws.Range("AF2").Value = 12.554
$ws.Range("AG2").Value = -353.266
$ws.Range("AH2").Value = 0.007161584607556116
$ws.Range("AI2").Value = 0.008119396652640646
$ws.Range("AJ2").Value = -0.2546714688597579
$ws.Range("AK2").Value = -0.2992881796478006
$ws.Range("AL2").Value = 0.986
$ws.Range("AM2").Value = 0.986
$ws.Range("AN2").Value = 0.04216996976822305
$ws.Range("AO2").Value = 286.0933062880325
$ws.Range("AP2").Value = -1.186650990930467
$ws.Range("AQ2").Value = 286.0933062880325
$ws.Range("D3").Value = 0.117
$ws.Range("E3").Value = 0.175
$ws.Range("F3").Value = 0.0285
$ws.Range("G3").Value = 0.1808338720103426
$ws.Range("H3").Value = 0.1808338720103426
$ws.Range("I3").Value = 0.1536845507433743
$ws.Range("J3").Value = 0.1312217194570136
$ws.Range("K3").Value = 80.5
$ws.Range("L3").Value = 0.1300904977375566
$ws.Range("M3").Value = 39.8
$ws.Range("N3").Value = 0.0397960203979602
$ws.Range("O3").Value = 0.4944099378881987
$ws.Range("P3").Value = 39.8
$ws.Range("Q3").Value = 0.0397960203979602
$ws.Range("R3").Value = 0.4944099378881987
$ws.Range("U3").Value = 125.4
$ws.Range("V3").Value = 0.1253874612538746
$ws.Range("W3").Value = 0.257847533632287
$ws.Range("X3").Value = 0.0481836194140207
$ws.Range("Y3").Value = 0.2096639142182663
$ws.Range("Z3").Value = 3.44223355788326
$ws.Range("AA3").Value = 0.4516958062380749
$ws.Range("AB3").Value = 0.04817774876382315
$ws.Range("AC3").Value = 0.4035180574742517
$ws.Range("AD3").Value = 0.23
$ws.Range("AE3").Value = 0
$ws.Range("AF3").Value = 0.23
$ws.Range("AG3").Value = -125.17
$ws.Range("AH3").Value = 0.0002299241250387372
$ws.Range("AI3").Value = 0.0006224122534029713
$ws.Range("AJ3").Value = -0.1430628736013167
$ws.Range("AK3").Value = -0.5127186335149306
$ws.Range("AN3").Value = 0.002224371373307543
$ws.Range("AP3").Value = -1.210541586073501
$ws.Range("D4").Value = 0.0788
$ws.Range("E4").Value = 0.116
$ws.Range("G4").Value = 0.1333687566418703
$ws.Range("H4").Value = 0.1333687566418703
$ws.Range("I4").Value = 0.1339001062699256
$ws.Range("J4").Value = 0.0933266232444166
$ws.Range("K4").Value = 122.6
$ws.Range("L4").Value = 0.09306209199939273
$ws.Range("M4").Value = 17
$ws.Range("N4").Value = 0.02613374327440431
$ws.Range("O4").Value = 0.1386623164763459
$ws.Range("P4").Value = 17
$ws.Range("Q4").Value = 0.02613374327440431
$ws.Range("R4").Value = 0.1386623164763459
$ws.Range("U4").Value = 231.4
$ws.Range("V4").Value = 0.3557263643351268
$ws.Range("W4").Value = 0.1577863577863578
$ws.Range("X4").Value = 0.04867645586533235
$ws.Range("Y4").Value = 0.1091099019210254
$ws.Range("Z4").Value = 2.225713803007265
$ws.Range("AA4").Value = 0.2077183535431567
$ws.Range("AB4").Value = 0.04824367569605216
$ws.Range("AC4").Value = 0.1594746778471045
$ws.Range("AD4").Value = 11
$ws.Range("AE4").Value = 0
$ws.Range("AF4").Value = 11
$ws.Range("AG4").Value = -220.4
$ws.Range("AH4").Value = 0.01662887377173091
$ws.Range("AI4").Value = 0.01097256857855362
$ws.Range("AJ4").Value = -0.5124389676819344
$ws.Range("AK4").Value = -0.2858254441706653
$ws.Range("AL4").Value = 0.925
$ws.Range("AM4").Value = 0.925
$ws.Range("AN4").Value = 0.06007646095030039
$ws.Range("AO4").Value = 190.7027027027027
$ws.Range("AP4").Value = -1.203713817586019
$ws.Range("AQ4").Value = 190.7027027027027
$ws.Range("D5").Value = -0.0443
$ws.Range("E5").Value = -0.163
$ws.Range("F5").Value = -0.243
$ws.Range("G5").Value = 0.2154531946508172
$ws.Range("H5").Value = 0.2154531946508172
$ws.Range("I5").Value = 0.1589895988112927
$ws.Range("J5").Value = 0.1376909921780818
$ws.Range("K5").Value = 6.74
$ws.Range("L5").Value = 0.100148588410104
$ws.Range("M5").Value = -0
$ws.Range("N5").Value = -0
$ws.Range("O5").Value = -0
$ws.Range("P5").Value = -0
$ws.Range("Q5").Value = -0
$ws.Range("R5").Value = -0
$ws.Range("U5").Value = 9.02
$ws.Range("V5").Value = 0.1108108108108108
$ws.Range("W5").Value = 0.05152905198776758
$ws.Range("X5").Value = 0.04864143590864292
$ws.Range("Y5").Value = 0.002887616079124661
$ws.Range("Z5").Value = 0.579622771509775
$ws.Range("AA5").Value = 0.07980883449819053
$ws.Range("AB5").Value = 0.04823906252611363
$ws.Range("AC5").Value = 0.03156977197207689
$ws.Range("AD5").Value = 1.28
$ws.Range("AE5").Value = 0
$ws.Range("AF5").Value = 1.28
$ws.Range("AG5").Value = -7.739999999999999
$ws.Range("AH5").Value = 0.01548137397194001
$ws.Range("AI5").Value = 0.007348719715237111
$ws.Range("AJ5").Value = -0.1050773825685582
$ws.Range("AK5").Value = -0.04686364737224509
$ws.Range("AL5").Value = 0.061
$ws.Range("AM5").Value = 0.061
$ws.Range("AN5").Value = 0.1142857142857143
$ws.Range("AO5").Value = 175.4098360655738
$ws.Range("AP5").Value = -0.6910714285714286
$ws.Range("AQ5").Value = 175.4098360655738
$ws.Range("A6").Value = "Malaysia"
$ws.Range("B6").Value = "Prevention Insurance.Com (OTCPK:PVNC)"
$ws.Range("C6").Value = "Insurance (General)"
$ws.Range("K6").Value = -0.112
$ws.Range("M6").Value = -0
$ws.Range("N6").Value = -0
$ws.Range("O6").Value = 0
$ws.Range("P6").Value = -0
$ws.Range("Q6").Value = -0
$ws.Range("R6").Value = 0
$ws.Range("S6").Value = 0
$ws.Range("U6").Value = 0
$ws.Range("V6").Value = 0
$ws.Range("W6").Value = 0.3227665706051874
$ws.Range("X6").Value = 0.04833140718427241
$ws.Range("Y6").Value = 0.274435163420915
$ws.Range("Z6").Value = 0
$ws.Range("AA6").Value = -55.99999999999995
$ws.Range("AB6").Value = 0.04822385452447246
$ws.Range("AC6").Value = -56.04822385452442
$ws.Range("AD6").Value = 0.044
$ws.Range("AE6").Value = 0
$ws.Range("AF6").Value = 0.044
$ws.Range("AG6").Value = 0.044
$ws.Range("AH6").Value = 0.005204636858291932
$ws.Range("AI6").Value = -1.222222222222222
$ws.Range("AJ6").Value = 0.005204636858291932
$ws.Range("AK6").Value = -1.222222222222222
$ws.Range("AL6").Value = 0
$ws.Range("AM6").Value = 0
